# The commit adds one more weekly price record for "Coliflor" (Cauliflower)
# at the "Macroferia Regional de Talca" market. The new record is inserted
# as row 65, which pushes every existing record from the old row 65 down
# through the old row 199 down by one row (to rows 66-200), growing the
# sheet's used range from A1:R199 to A1:R200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 65; Excel shifts rows 65..199 down to
# 66..200 automatically (carrying their formatting, including the date
# style on column D).
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record's data.
$ws.Cells.Item(65, 1).Value2  = 5
$ws.Cells.Item(65, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value2  = "Maule"
$ws.Cells.Item(65, 4).Value2  = 44614
$ws.Cells.Item(65, 5).Value2  = 7
$ws.Cells.Item(65, 6).Value2  = 100112008
$ws.Cells.Item(65, 7).Value2  = "Coliflor"
$ws.Cells.Item(65, 8).Value2  = "Sin especificar"
$ws.Cells.Item(65, 9).Value2  = "Segunda"
$ws.Cells.Item(65, 10).Value2 = 2000
$ws.Cells.Item(65, 11).Value2 = 1200
$ws.Cells.Item(65, 12).Value2 = 1200
$ws.Cells.Item(65, 13).Value2 = 1200
$ws.Cells.Item(65, 14).Value2 = "$/unidad"
$ws.Cells.Item(65, 15).Value2 = "Región del Maule"
$ws.Cells.Item(65, 16).Value2 = 1200
$ws.Cells.Item(65, 17).Value2 = 1
$ws.Cells.Item(65, 18).Value2 = "Hortaliza"
